$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A63").Value = "Francesca Sartori"
$ws.Range("B63").Value = "Elia Barozzi | I Magnifici"
$ws.Range("C63").Value = "Edoardo Pomarolli | Modium"
$ws.Range("D63").Value = "Giacomo Gasparini | MAI UNA GIOIA"
$ws.Range("E63").Value = "Leonardo  Parisi  | MediaserT"
$ws.Range("F63").Value = "Daniele Dalbosco | SdrumALA"
